# Updates hotel review data: hotel_info row 2 review counts, and 3 new review_info rows
$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# --- hotel_info: fill in review counters on row 2 ---
$wsHotel = $wb.Worksheets.Item("hotel_info")
Set-TextCell $wsHotel 'G2' '7'
Set-TextCell $wsHotel 'H2' '3'
Set-TextCell $wsHotel 'I2' '7'

# --- review_info: append 3 new review rows ---
$wsReview = $wb.Worksheets.Item("review_info")

# Row 2
$wsReview.Range("A2").Value = 66692
$wsReview.Range("D2").Value = 1
Set-TextCell $wsReview 'E2' '08/06/2018'
Set-TextCell $wsReview 'F2' 'https://www.tripadvisor.com/ShowUserReviews-g55609-d12944276-r577975854-Comfort_Suites-Channelview_Texas.html'
Set-TextCell $wsReview 'G2' '55609'
Set-TextCell $wsReview 'H2' '12944276'
Set-TextCell $wsReview 'I2' '577975854'
Set-TextCell $wsReview 'J2' '05/04/2018'
Set-TextCell $wsReview 'K2' 'Clean and comfortable'
Set-TextCell $wsReview 'L2' 'My experience with this property was great. From check in to check out this site is worth the stay. I found the room to be clean and felt comfortable. It was easy to find just off the highway.The price was just right.'
$wsReview.Range("M2").Value = 5
Set-TextCell $wsReview 'N2' 'May 2018'
Set-TextCell $wsReview 'O2' ' traveled on business'
$wsReview.Range("R2").Value = 5
$wsReview.Range("S2").Value = 5
$wsReview.Range("U2").Value = 5
$wsReview.Range("V2").Value = 0
Set-TextCell $wsReview 'Y2' 'My experience with this property was great. From check in to check out this site is worth the stay. I found the room to be clean and felt comfortable. It was easy to find just off the highway.The price was just right.'

# Row 3
$wsReview.Range("A3").Value = 66692
$wsReview.Range("D3").Value = 2
Set-TextCell $wsReview 'E3' '08/06/2018'
Set-TextCell $wsReview 'F3' 'https://www.tripadvisor.com/ShowUserReviews-g55609-d12944276-r567371609-Comfort_Suites-Channelview_Texas.html'
Set-TextCell $wsReview 'G3' '55609'
Set-TextCell $wsReview 'H3' '12944276'
Set-TextCell $wsReview 'I3' '567371609'
Set-TextCell $wsReview 'J3' '03/18/2018'
Set-TextCell $wsReview 'K3' 'A number 1'
Set-TextCell $wsReview 'L3' 'This is a new property, super clean, great rooms. The staff was awesome, check in was quick and easy. It was so quiet and comfortable my wife and I almost overslept. There are also numerous places to eat with a 10 minute drive.'
$wsReview.Range("M3").Value = 5
Set-TextCell $wsReview 'N3' 'March 2018'
Set-TextCell $wsReview 'O3' ' traveled with family'
$wsReview.Range("V3").Value = 0
Set-TextCell $wsReview 'Y3' 'This is a new property, super clean, great rooms. The staff was awesome, check in was quick and easy. It was so quiet and comfortable my wife and I almost overslept. There are also numerous places to eat with a 10 minute drive.'

# Row 4
$wsReview.Range("A4").Value = 66692
$wsReview.Range("D4").Value = 3
Set-TextCell $wsReview 'E4' '08/06/2018'
Set-TextCell $wsReview 'F4' 'https://www.tripadvisor.com/ShowUserReviews-g55609-d12944276-r549981689-Comfort_Suites-Channelview_Texas.html'
Set-TextCell $wsReview 'G4' '55609'
Set-TextCell $wsReview 'H4' '12944276'
Set-TextCell $wsReview 'I4' '549981689'
Set-TextCell $wsReview 'J4' '12/29/2017'
Set-TextCell $wsReview 'K4' 'Great Staff Here! Awesome Value!'
Set-TextCell $wsReview 'L4' 'This hotel is very new and is super clean. The room was great! Spacious clean & modern. it felt more like a room valued at $200.  I was very impressed with the team at this hotel. They were great communicators and seemed genuinely interested in your satisfaction. I will absolutely stay here again when in the Houston area. Overall it''s one the best hotels I have ever stayed in.'
$wsReview.Range("M4").Value = 5
Set-TextCell $wsReview 'N4' 'December 2017'
Set-TextCell $wsReview 'O4' ' traveled on business'
$wsReview.Range("U4").Value = 5
$wsReview.Range("Q4").Value = 5
$wsReview.Range("V4").Value = 0
Set-TextCell $wsReview 'Y4' 'This hotel is very new and is super clean. The room was great! Spacious clean & modern. it felt more like a room valued at $200.  I was very impressed with the team at this hotel. They were great communicators and seemed genuinely interested in your satisfaction. I will absolutely stay here again when in the Houston area. Overall it''s one the best hotels I have ever stayed in.'

